$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.60121762752533
$ws.Range("B1").Value = 1.967608690261841
$ws.Range("C1").Value = 2.135896444320679
$ws.Range("D1").Value = 2.460945129394531
$ws.Range("E1").Value = 3.271326065063477
